# Update RQ134 variable definitions for September eval
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Levels" text for the Delegator_grp variable (row 4, column C)
$ws.Cells.Item(4, 3).Value = "P2E_June_2025, P2E_July_2025"

# Update the "Levels" text for the Attribute variable (row 4, column H)
$ws.Cells.Item(4, 8).Value = "Affiliation Focus (AF), Merit Focus (MF), Personal Safety (PS), Search vs. Stay (SS), Affiliation Focus / Merit Focus (AF-MF)"

# Re-fit row heights to account for the new wrapped text lengths
$ws.Rows.Item(1).RowHeight = 31.5
$ws.Rows.Item(2).RowHeight = 44.25
$ws.Rows.Item(3).RowHeight = 171.75
$ws.Rows.Item(4).RowHeight = 159
